# V31: Engraved bars 161-168; corrected mistakes from last week
#
# Appends 11 new rows to the single errata table, one per new proofreading
# issue (Instrument | Bar | Issue | Answer). The "Answer" column is left
# blank, matching the existing rows' pattern.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRows = @(
    @("Viola", "151", "Treble clef from bar 153 moved forward to avoid clef change over tied notes. "),
    @("Cello", "151", "Second note head on G added to indicate unison double stop. "),
    @("All", "151-2", "4/4 time signature added"),
    @("All", "153", "Do the glissando gestures lead to the first notes in bar 154, or do they end on an indeterminate pitch before the attack on bar 154? "),
    @("Viola", "154", "Would you like the viola to start at ff? All other parts start at ff"),
    @("Violins and Cello", "154", "Would you like the highlighted notes to be marked tenuto and with accents? "),
    @("All", "157", "Would you like an fp indication on every note? Or is this bar to be played piano throughout, with just a forte attack at the beginning? "),
    @("Violin II", "157-8", "Would you like this passage to be taken under one bow, or should the bow direction be changed for each new accent? "),
    @("Violin I", "160", "Which note would you like dotted? "),
    @("Viola", "162", "Which note would you like dotted? "),
    @("Violin II", "166-7", "Would you like the lower voice to be accented as well? ")
)

foreach ($rowData in $newRows) {
    $newRow = $t.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $rowData[0]
    $newRow.Cells.Item(2).Range.Text = $rowData[1]
    $newRow.Cells.Item(3).Range.Text = $rowData[2]
    # Cell 4 ("Answer") is left blank, same as every preceding row.
}
